$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right after "总计" (shifts all the
#    quarterly sheets down by one tab position).
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$placeholder = $wb.Worksheets.Add($null, $afterSheet)
$placeholder.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add a new 2022-Q4 row at the top of
#    the data block and shift the rest down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Extend the A-column "index" style (bold/centered/bordered) down into the
# brand new row 8 before writing values, by cloning the format already used
# by row 7.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q4", 40, 6.01),
    @(1, "2022-Q3", 14, 2.2),
    @(2, "2022-Q2", 3, 0.03),
    @(3, "2021-Q3", 1, 0.19),
    @(4, "2021-Q2", 5, 0.29),
    @(5, "2021-Q1", 4, 0.29),
    @(6, "2020-Q4", 4, 1.95)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 3. Populate the new "2022-Q4" sheet with the fund holdings table.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2022-Q4")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q4Rows = @(
    @(0, '010420', '民生加银成长优选股票', '33.68', '82.33', '3.91', '1.3169', 8),
    @(1, '000136', '民生加银策略精选混合A', '19.52', '81.96', '4.24', '0.8276', 5),
    @(2, '070021', '嘉实主题新动力混合', '17.08', '93.12', '4.23', '0.7225', 5),
    @(3, '000985', '嘉实逆向策略股票', '11.36', '93.39', '4.08', '0.4635', 5),
    @(4, '003378', '泰康策略优选灵活配置混合', '13.78', '83.09', '2.99', '0.4120', 8),
    @(5, '009659', '民生加银新动能一年定期开放混合A', '8.95', '82.72', '4.42', '0.3956', 5),
    @(6, '012466', '嘉实策略精选混合A', '9.59', '93.89', '4.10', '0.3932', 4),
    @(7, '010116', '民生加银新兴产业混合A', '6.99', '85.62', '4.75', '0.3320', 5),
    @(8, '009660', '民生加银新动能一年定期开放混合C', '3.45', '82.72', '4.42', '0.1525', 5),
    @(9, '006058', '民生加银新兴成长混合', '3.44', '86.91', '4.15', '0.1428', 4),
    @(10, '501200', '民生加银科技创新 3 年封闭混合', '2.45', '86.66', '4.04', '0.0990', 6),
    @(11, '005310', '广发电子信息传媒产业精选股票A', '1.72', '88.90', '5.03', '0.0865', 2),
    @(12, '001170', '泰达宏利复兴伟业灵活配置混合', '1.64', '91.15', '4.87', '0.0799', 6),
    @(13, '007853', '华商计算机行业量化股票A', '2.18', '92.32', '3.52', '0.0767', 4),
    @(14, '014307', '嘉实多元动力混合A', '1.56', '93.22', '4.08', '0.0636', 4),
    @(15, '004671', '中融核心成长灵活配置混合', '1.14', '65.86', '5.01', '0.0571', 2),
    @(16, '010009', '中融成长优选混合C', '1.05', '60.81', '4.65', '0.0488', 1),
    @(17, '012467', '嘉实策略精选混合C', '1.01', '93.89', '4.10', '0.0414', 4),
    @(18, '016029', '湘财成长优选一年持有期混合A', '1.70', '89.45', '2.37', '0.0403', 10),
    @(19, '010117', '民生加银新兴产业混合C', '0.76', '85.62', '4.75', '0.0361', 5),
    @(20, '013262', '西部利得个股精选股票C', '1.25', '91.27', '2.67', '0.0334', 5),
    @(21, '673090', '西部利得个股精选股票A', '1.22', '91.27', '2.67', '0.0326', 5),
    @(22, '010008', '中融成长优选混合A', '0.57', '60.81', '4.65', '0.0265', 1),
    @(23, '000966', '中邮核心科技创新灵活配置混合', '0.96', '87.45', '2.28', '0.0219', 10),
    @(24, '008336', '宝盈祥裕增强回报混合A', '0.69', '36.55', '2.28', '0.0157', 9),
    @(25, '014308', '嘉实多元动力混合C', '0.33', '93.22', '4.08', '0.0135', 4),
    @(26, '350005', '天治中国制造2025灵活配置混合', '0.24', '94.00', '5.05', '0.0121', 6),
    @(27, '010236', '广发电子信息传媒产业精选股票C', '0.24', '88.90', '5.03', '0.0121', 2),
    @(28, '519127', '浦银安盛盛世精选灵活配置混合A', '1.17', '23.43', '0.91', '0.0106', 5),
    @(29, '000994', '建信睿盈灵活配置混合A', '0.33', '90.01', '3.14', '0.0104', 10),
    @(30, '519177', '浦银安盛盛世精选灵活配置混合C', '0.76', '23.43', '0.91', '0.0069', 5),
    @(31, '009709', '民生加银策略精选混合C', '0.16', '81.96', '4.24', '0.0068', 5),
    @(32, '000995', '建信睿盈灵活配置混合C', '0.16', '90.01', '3.14', '0.0050', 10),
    @(33, '519175', '浦银安盛经济带崛起灵活配置混合', '0.39', '23.04', '1.12', '0.0044', 2),
    @(34, '016030', '湘财成长优选一年持有期混合C', '0.15', '89.45', '2.37', '0.0036', 10),
    @(35, '008337', '宝盈祥裕增强回报混合C', '0.08', '36.55', '2.28', '0.0018', 9),
    @(36, '001530', '万家瑞富灵活配置混合A', '0.21', '23.67', '0.74', '0.0016', 5),
    @(37, '008162', '浦银安盛经济带崛起灵活配置混合C', '0.12', '23.04', '1.12', '0.0013', 2),
    @(38, '012007', '万家瑞富灵活配置混合C', '0.11', '23.67', '0.74', '0.0008', 5),
    @(39, '017628', '华商计算机行业量化股票C', '0.00', '92.32', '3.52', '0', 4)
)

# Mark D:G as text first so the numeric-looking strings ("33.68", "82.33", ...)
# are stored verbatim instead of being coerced into numbers, then strip the
# format override back off so the cells keep the workbook's default style.
$textRange = $q4.Range("D2:G41")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $q4Rows.Length; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[7]
}

$textRange.ClearFormats()

# Row 41's "持有市值" is a genuine number (0), not text like the other rows.
for ($i = 0; $i -lt $q4Rows.Length; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    if ($i -eq 39) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = $row[6]
    }
}

# Clone header (B1:H1) and index-column (A2:A41) formatting from the sheet
# that already carries it (the old "2022-Q3" sheet, still at its original
# name at this point in the script).
$styleSource = $wb.Worksheets.Item("2022-Q3")
$styleSource.Range("B1:H1").Copy()
$q4b = $wb.Worksheets.Item("2022-Q4")
$q4b.Range("B1:H1").PasteSpecial(-4122)

$styleSource2 = $wb.Worksheets.Item("2022-Q3")
$styleSource2.Range("A2:A15").Copy()
$q4c = $wb.Worksheets.Item("2022-Q4")
$q4c.Range("A2:A41").PasteSpecial(-4122)
